$d = $word.ActiveDocument

# Sequentially replace text occurrences, always searching forward from the
# end of the previous match so that repeated substrings (like "[pause]")
# are resolved in document order rather than Find.Execute's "replace all"
# or "always match the first occurrence" behavior.
$script:pos = 0

function Replace-Next($find, $replace) {
    $searchRange = $d.Range($script:pos, $d.Content.End)
    $found = $searchRange.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 1)
    if ($found) {
        $script:pos = $searchRange.End
    }
}

Replace-Next "neck, " "گردن، "
Replace-Next "[pause] " "[مکث] "
Replace-Next "shoulders. " "شانه ها. "
Replace-Next "Noticing any sensations, movements, any places of holding. " "به هر گونه احساس، حرکت یا ناحیه که تنش دارد توجه کنید. "

Replace-Next "Continue to scan your body, moving down the arms and hands " "چک کردن را ادامه دهید، از بازو ها ‌و دست ها عبوز کنید "
Replace-Next "[pause] " "[مکث] "
Replace-Next "and into your fingers. " "تا نوک انگشتان تان. "
Replace-Next "Sense the back of your body" "پشت بدن تانرا «کمر» تان را احساس کنید"
Replace-Next "[pause]" "[مکث]"
Replace-Next " and your lower back. " " وقسمت پایین کمر تان. "

Replace-Next "Feeling the contact of your body with the chair if you are seated. " "اگر نشسته هستید، تماس بدن تانرا با چوکی احساس کنید. "
Replace-Next "Now, sensing your upper legs, " "حال توجه و احساس تان به ران ها، "
Replace-Next "[pause]" "[مکث]"
Replace-Next "your lower legs, " "ساق ‌پاه ها، "
Replace-Next "[pause] " "[مکث] "
Replace-Next "and the feet. " "و پاها منتقل کنید. "

Replace-Next "Take a full, deep breath in and a long breath out." "یک نفس کامل و عمیق بکشید و آرامی بیرون دهید."
Replace-Next "Finish this pause by opening your eyes. " "این مکث را با باز کردن چشمان تان به پایان برسانید. "

Replace-Next "Thank you for taking a moment to pause with us. It's time for today's lesson. " "تشکر ازینکه با ما لحظه ای مکث کردید. اکنون وقت درس امروز است. "

Replace-Next "Today's pause is called loving kindness meditation." "مکث امروز مدیتیشن محبت‌آمیز نام دارد."

Replace-Next "Find a comfortable sitting position, your feet flat on the floor, your hands resting in your lap." "در یک حالت نشسته‌ی راحت قرار بگیرید، پاها را صاف روی زمین بگذارید و دست‌ها را روی زانو یا در دامنتان بگذارید."
Replace-Next "Close your eyes, if you are comfortable to do so, or allow your eyes and eyelids to relax and your gaze to become soft." "چشمان خود را ببندید، اگر با آن راحت هستید، یا بگذارید پلک‌ها و نگاهتان نرم و آرام شود."
Replace-Next "Ask yourself, ""What is my experience in this moment?"" " "از خود بپرسید: «همین حالا چه احساسی دارم؟» "

Replace-Next "Notice what thoughts you are experiencing. " "به افکاری که در ذهن تان هستند توجه کنید. "
Replace-Next "Notice how you feel emotionally. " "احساسات تان را شناسایی کنید. "
Replace-Next "Notice how your body feels. Notice any discomfort or tension." "بدن تان را احساس کنید. به تنش و ناراحتی که است را توجه کنید."

Replace-Next "Connect to your heart in a kind and gentle way. You may want to place one hand on your heart or chest. " "با مهربانی و نرمی با قلب خود ارتباط برقرار کنید. می‌توانید یک دست‌ تان را روی قلب یا سینه ‌تان بگذارید. "
Replace-Next "You can then say the following words silently to yourself: " "سپس می‌توانید این جملات را در دل خود آهسته تکرار کنید: "
